$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "Adaptive Solar Facade" header (column D) to reflect
# the first cluster configuration, and add a new column E for the second
# cluster configuration (new Geneva combination).
$ws.Range("D1").Value2 = "Adaptive Solar Facade (2 Clusters, 3*6 angles)"
$ws.Range("E1").Value2 = "Adaptive Solar Facade (1 Cluster, 3*12 angles)"

# New Geneva results for the "Adaptive Solar Facade (1 Cluster, 3*12 angles)" column
$ws.Range("E16").Value2 = 442.11
$ws.Range("E17").Value2 = 504.38
$ws.Range("E18").Value2 = 432.6
$ws.Range("E19").Value2 = 1379.09

# Match the bold "Total" style already used by D19/C19/B19
$ws.Range("E19").Font.Bold = $true

# Widen columns D and E to fit the longer headers
$ws.Columns.Item(4).ColumnWidth = 40.83
$ws.Columns.Item(5).ColumnWidth = 39.83

# Update the active selection to reflect where the author ended up working
$ws.Range("E20").Select() | Out-Null
